$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.283484101295471
$ws.Range("B1").Value = 2.345788717269897
$ws.Range("C1").Value = 2.992778778076172
$ws.Range("D1").Value = 3.450932502746582
$ws.Range("E1").Value = 1.568946957588196
